$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'243.86"
$ws.Range("D3").Formula = "'23.96"
$ws.Range("D4").Formula = "'5.139"
$ws.Range("D5").Formula = "'0.05745"
$ws.Range("D6").Formula = "'6.484"
$ws.Range("D7").Formula = "'3.144"
$ws.Range("D8").Formula = "'0.8107"
$ws.Range("D9").Formula = "'0.8349"
$ws.Range("D10").Formula = "'0.1339"
$ws.Range("D11").Formula = "'0.06946"
$ws.Range("D12").Formula = "'0.03125"
$ws.Range("D13").Formula = "'0.02846"
$ws.Range("D14").Formula = "'0.09366"
$ws.Range("D15").Formula = "'3.749"
$ws.Range("D16").Formula = "'0.001514"
$ws.Range("D17").Formula = "'0.04659"
$ws.Range("D18").Formula = "'0.0005965"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Formula = "'0.006145"
$ws.Range("D20").Formula = "'0.001237"
$ws.Range("D21").Formula = "'0.004276"
$ws.Range("D22").Formula = "'0.00008709"
$ws.Range("D24").Formula = "'2.082"
$ws.Range("D25").Formula = "'0.3172"
$ws.Range("D26").Formula = "'0.1338"
$ws.Range("D28").Formula = "'0.0002331"
$ws.Range("D40").Formula = "'0.03616"
$ws.Range("D41").Formula = "'0.006306"
$ws.Range("D42").Formula = "'0.1049"
$ws.Range("D43").Formula = "'0.002934"
$ws.Range("D44").Formula = "'0.007360"
$ws.Range("D45").Formula = "'0.00005300"
$ws.Range("D47").Formula = "'0.2672"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Formula = "'0.002277"
$ws.Range("D49").Formula = "'0.00002102"
$ws.Range("D50").Formula = "'0.0002001"
